$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 4 and 5 (Sending cluster = MuSCs rows using old TPM data);
# their replacement data will now live in rows 2-3 using updated TPM-derived values.
$ws.Rows("4:5").Delete()

# Row 2: MuSCs -> ECs (Il10/Il10ra), refreshed with new TPM-based values
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Il10"
$ws.Range("C2").Value = "Il10ra"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08735766666666667
$ws.Range("H2").Value = 0.262073
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1280766666666666
$ws.Range("N2").Value = 0.38423
$ws.Range("O2").Value = 0.9809242740654885
$ws.Range("P2").Value = 0.9809242740654887
$ws.Range("Q2").Value = 0.01118847875444444
$ws.Range("R2").Value = 0.10069630879
$ws.Range("S2").Value = 0.9809242740654885
$ws.Range("T2").Value = 0.9809242740654887

# Row 3: MuSCs -> MuSCs (Il10/Il10ra), refreshed with new TPM-based values
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Il10"
$ws.Range("C3").Value = "Il10ra"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08735766666666667
$ws.Range("H3").Value = 0.262073
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.002490666666666667
$ws.Range("N3").Value = 0.007472
$ws.Range("O3").Value = 0.01907572593451144
$ws.Range("P3").Value = 0.01907572593451144
$ws.Range("Q3").Value = 0.0002175788284444444
$ws.Range("R3").Value = 0.001958209456
$ws.Range("S3").Value = 0.01907572593451144
$ws.Range("T3").Value = 0.01907572593451144
